$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "93.387.60"
$ws.Range("E2").Value = "  +1.33%  "

# Row 3
$ws.Range("D3").Value = "3.460.48"
$ws.Range("E3").Value = "  +4.18%  "

# Row 4
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").Value = "'234.68"
$ws.Range("E5").Value = "  +3.33%  "

# Row 6
$ws.Range("D6").Value = "'625.44"
$ws.Range("E6").Value = "  +0.96%  "

# Row 7
$ws.Range("E7").Value = "  +7.78%  "

# Row 8
$ws.Range("D8").Value = "'0.390"
$ws.Range("E8").Value = "  +3.77%  "

# Row 9
$ws.Range("E9").Value = "  +0.09%  "

# Row 10
$ws.Range("E10").Value = "  +11.49%  "

# Row 11
$ws.Range("D11").Value = "3.455.65"
$ws.Range("E11").Value = "  +4.04%  "

# Row 12
$ws.Range("D12").Value = "'43.37"
$ws.Range("E12").Value = "  +9.81%  "

# Row 13
$ws.Range("E13").Value = "  +5.37%  "

# Row 14
$ws.Range("E14").Value = "  +7.18%  "

# Row 15
$ws.Range("D15").Value = "4.113.31"
$ws.Range("E15").Value = "  +4.46%  "

# Row 16
$ws.Range("D16").Value = "93.214.34"
$ws.Range("E16").Value = "  +1.36%  "

# Row 17
$ws.Range("D17").Value = "'0.0000248"
$ws.Range("E17").Value = "  +3.53%  "

# Row 18
$ws.Range("D18").Value = "'8.30"
$ws.Range("E18").Value = "  +6.03%  "

# Row 19
$ws.Range("D19").Value = "3.456.43"
$ws.Range("E19").Value = "  +4.41%  "

# Row 20
$ws.Range("D20").Value = "'18.07"
$ws.Range("E20").Value = "  +9.24%  "

# Row 21
$ws.Range("D21").Value = "'11.77"
$ws.Range("E21").Value = "  +9.87%  "

# Row 22
$ws.Range("D22").Value = "'0.500"
$ws.Range("E22").Value = "  +14.73%  "

# Row 23
$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").Value = "'3.38"
$ws.Range("E23").Value = "  +10.99%  "

# Row 24
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "'502.85"
$ws.Range("E24").Value = "  +3.99%  "

# Row 25
$ws.Range("D25").Value = "'6.79"
$ws.Range("E25").Value = "  +12.75%  "

# Row 26
$ws.Range("D26").Value = "'0.0000183"
$ws.Range("E26").Value = "  +2.60%  "

# Row 27
$ws.Range("D27").Value = "'94.83"
$ws.Range("E27").Value = "  +7.32%  "

# Row 28
$ws.Range("D28").Value = "'12.13"
$ws.Range("E28").Value = "  +8.08%  "

# Row 29
$ws.Range("D29").Value = "3.643.77"
$ws.Range("E29").Value = "  +4.70%  "

# Row 30
$ws.Range("E30").Value = "  +10.27%  "

# Row 31
$ws.Range("D31").Value = "'11.32"
$ws.Range("E31").Value = "  +2.46%  "

# Row 32
$ws.Range("B32").Value = "Dai"
$ws.Range("C32").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.04%  "

# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.140"
$ws.Range("E33").Value = "  +7.88%  "

# Row 34
$ws.Range("D34").Value = "'0.992"
$ws.Range("E34").Value = "  -0.45%  "

# Row 35
$ws.Range("E35").Value = "  +7.52%  "

# Row 36
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "'29.47"
$ws.Range("E36").Value = "  +4.98%  "

# Row 37
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").Value = "'0.555"
$ws.Range("E37").Value = "  +8.32%  "

# Row 38
$ws.Range("D38").Value = "'572.71"
$ws.Range("E38").Value = "  +10.61%  "

# Row 39
$ws.Range("E39").Value = "  +6.37%  "

# Row 40
$ws.Range("D40").Value = "'7.52"
$ws.Range("E40").Value = "  +3.48%  "

# Row 42
$ws.Range("D42").Value = "'0.913"
$ws.Range("E42").Value = "  +6.46%  "

# Row 43
$ws.Range("E43").Value = "  +2.36%  "

# Row 44
$ws.Range("D44").Value = "'23.73"
$ws.Range("E44").Value = "  -0.96%  "

# Row 45
$ws.Range("D45").Value = "'0.0421"
$ws.Range("E45").Value = "  +9.55%  "

# Row 46
$ws.Range("E46").Value = "  +3.75%  "

# Row 47
$ws.Range("D47").Value = "'5.49"
$ws.Range("E47").Value = "  +3.48%  "

# Row 48
$ws.Range("E48").Value = "  -0.87%  "

# Row 49
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "'53.22"
$ws.Range("E49").Value = "  +2.59%  "

# Row 50
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'2.12"
$ws.Range("E50").Value = "  +0.42%  "

# Row 51
$ws.Range("D51").Value = "'8.12"
$ws.Range("E51").Value = "  +5.12%  "
